{"js": "const body = context.document.body;\n\n// Update the date paragraph (first paragraph in the document).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nconst datePara = paragraphs.items[0];\ndatePara.insertText(\"2025-03-17 Monday\", Word.InsertLocation.replace);\n\n// Update every cell in the practice table (20 rows x 5 columns) in one shot,\n// preserving each cell's existing formatting (font/size/paragraph alignment).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\ntable.values = [[\"86-71=15\", \"38-31=7\", \"60-17=43\", \"22+28=50\", \"74-40=34\"], [\"61-21=40\", \"98-84=14\", \"18+62=80\", \"84+12=96\", \"59+17=76\"], [\"28+11=39\", \"78-62=16\", \"73-67=6\", \"87+2=89\", \"91-67=24\"], [\"97-89=8\", \"46+30=76\", \"58-51=7\", \"85+8=93\", \"3+75=78\"], [\"66-29=37\", \"53-39=14\", \"13+85=98\", \"34+38=72\", \"5+75=80\"], [\"17+25=42\", \"39+3=42\", \"77-23=54\", \"33+48=81\", \"98-46=52\"], [\"96-35=61\", \"56+29=85\", \"71-64=7\", \"42-24=18\", \"14+61=75\"], [\"96-69=27\", \"49-0=49\", \"54+20=74\", \"18+28=46\", \"39+35=74\"], [\"73-70=3\", \"39-10=29\", \"72-7=65\", \"75+10=85\", \"37+5=42\"], [\"48-8=40\", \"93-18=75\", \"37+8=45\", \"91-14=77\", \"96-95=1\"], [\"92-89=3\", \"71+21=92\", \"34+55=89\", \"98-45=53\", \"93-52=41\"], [\"41+24=65\", \"89-67=22\", \"91-90=1\", \"11+59=70\", \"18+56=74\"], [\"31-11=20\", \"10+62=72\", \"6+75=81\", \"69-53=16\", \"42+38=80\"], [\"80-20=60\", \"1+51=52\", \"72+15=87\", \"95-9=86\", \"16+63=79\"], [\"3+73=76\", \"51+48=99\", \"21+3=24\", \"66-27=39\", \"69-28=41\"], [\"60-9=51\", \"60-9=51\", \"99-14=85\", \"34+8=42\", \"71-71=0\"], [\"38-7=31\", \"52-30=22\", \"58-14=44\", \"12+1=13\", \"49-21=28\"], [\"57+22=79\", \"1+2=3\", \"36+38=74\", \"7+42=49\", \"84-12=72\"], [\"30-23=7\", \"98-74=24\", \"52-35=17\", \"6+23=29\", \"84-13=71\"], [\"18+50=68\", \"42+34=76\", \"5+60=65\", \"50+2=52\", \"61+32=93\"]];\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date paragraph (first paragraph in the document).\n$d.Paragraphs.Item(1).Range.Text = \"2025-03-17 Monday\"\n\n# Update every cell in the practice table (20 rows x 5 columns), preserving\n# each cell's existing formatting (font/size/paragraph alignment).\n$table = $d.Tables.Item(1)\n$values = @(\n    @(\"86-71=15\",\"38-31=7\",\"60-17=43\",\"22+28=50\",\"74-40=34\"),\n    @(\"61-21=40\",\"98-84=14\",\"18+62=80\",\"84+12=96\",\"59+17=76\"),\n    @(\"28+11=39\",\"78-62=16\",\"73-67=6\",\"87+2=89\",\"91-67=24\"),\n    @(\"97-89=8\",\"46+30=76\",\"58-51=7\",\"85+8=93\",\"3+75=78\"),\n    @(\"66-29=37\",\"53-39=14\",\"13+85=98\",\"34+38=72\",\"5+75=80\"),\n    @(\"17+25=42\",\"39+3=42\",\"77-23=54\",\"33+48=81\",\"98-46=52\"),\n    @(\"96-35=61\",\"56+29=85\",\"71-64=7\",\"42-24=18\",\"14+61=75\"),\n    @(\"96-69=27\",\"49-0=49\",\"54+20=74\",\"18+28=46\",\"39+35=74\"),\n    @(\"73-70=3\",\"39-10=29\",\"72-7=65\",\"75+10=85\",\"37+5=42\"),\n    @(\"48-8=40\",\"93-18=75\",\"37+8=45\",\"91-14=77\",\"96-95=1\"),\n    @(\"92-89=3\",\"71+21=92\",\"34+55=89\",\"98-45=53\",\"93-52=41\"),\n    @(\"41+24=65\",\"89-67=22\",\"91-90=1\",\"11+59=70\",\"18+56=74\"),\n    @(\"31-11=20\",\"10+62=72\",\"6+75=81\",\"69-53=16\",\"42+38=80\"),\n    @(\"80-20=60\",\"1+51=52\",\"72+15=87\",\"95-9=86\",\"16+63=79\"),\n    @(\"3+73=76\",\"51+48=99\",\"21+3=24\",\"66-27=39\",\"69-28=41\"),\n    @(\"60-9=51\",\"60-9=51\",\"99-14=85\",\"34+8=42\",\"71-71=0\"),\n    @(\"38-7=31\",\"52-30=22\",\"58-14=44\",\"12+1=13\",\"49-21=28\"),\n    @(\"57+22=79\",\"1+2=3\",\"36+38=74\",\"7+42=49\",\"84-12=72\"),\n    @(\"30-23=7\",\"98-74=24\",\"52-35=17\",\"6+23=29\",\"84-13=71\"),\n    @(\"18+50=68\",\"42+34=76\",\"5+60=65\",\"50+2=52\",\"61+32=93\"),\n)\n\nfor ($r = 1; $r -le $values.Count; $r++) {\n    $row = $values[$r - 1]\n    for ($c = 1; $c -le $row.Count; $c++) {\n        $table.Cell($r, $c).Range.Text = $row[$c - 1]\n    }\n}\n"}
